# ER-Diagram.xlsx: add a "ForgotPassword" table (Fpid PK, UserID FK, Date, Code)
# to the top of the ER diagram sheet, above the existing "Users" table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = "ForgotPassword"
$ws.Range("B4").Value = "Fpid INT (PK)"
$ws.Range("B5").Value = "UserID INT (FK)"
$ws.Range("B6").Value = "Date DATE"
$ws.Range("B7").Value = "Code NvarChar(8)"

# Match the author's final selection/scroll position.
$ws.Range("E3").Select()
